$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set G3 (Invalid) to 1
$ws.Range("G3").Value = 1

# Set H3:H18 (Absent) to 1 for each date row
$ws.Range("H3:H18").Value = 1
